$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from row 6 so the new row matches the existing table style
$ws.Range("A6:C6").Copy()
$ws.Range("A7:C7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New row 7 data: date 2021-09-01 (serial 44440), Toallitas 14.34, Tampones 16.38
$ws.Range("A7").Value2 = 44440
$ws.Range("B7").Value2 = 14.34
$ws.Range("C7").Value2 = 16.38

# Update the active selection to reflect the edit location
$ws.Range("C12").Select()
